# Generate Report for Handoff
# Replaces the stale "d175820d...md" handoff file references with the new
# "821d87f1...md" ones, drops the failed-transform dependency row
# ("974e5648...md" / "Handoff transform failed"), and refreshes the
# handoff xlf file names / timestamps for both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$oldGuid = "d175820d-95d4-48af-8961-accbf0f466f9"
$newGuid = "821d87f1-8817-4972-a132-099f001c662d"
$oldHash = "a6a884df89f2d9ef69d097566ac56bf43d933e15"
$newHash = "749c0e0d1bd1a1085f7599b9311e01b281d2ef6f"

$mdBase = "https://github.com/OpenLocalizationTest/oltest/blob/1f18b90e9bfc546c808228ed0f2f27c43d0893a4/e2e/"
$cfgUrl = "https://github.com/OpenLocalizationTest/oltest/blob/1f18b90e9bfc546c808228ed0f2f27c43d0893a4/.localization-config"
$zhXlfBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7ffde657cd55425b88d4e08c17956bb94921de06/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/"
$deXlfBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0ab417a062f1a726b9adda254a020b8503fc2ef4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/"

$newMdName = $newGuid + ".md"
$newZhXlfName = $newGuid + "." + $newHash + ".zh-cn.xlf"
$newDeXlfName = $newGuid + "." + $newHash + ".de-de.xlf"

$newZhDatetime = "2016-01-27 08:35:06"
$newDeDatetime = "2016-01-27 08:35:25"

# ---------------------------------------------------------------
# Sheet "Overview" : drop row 3 (974e5648.../Handoff transform failed),
# row 4 (.localization-config) shifts up to become row 3. Update the
# file name on (old) row 2.
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Hyperlinks.Delete()
$wsOverview.Range("A3").Hyperlinks.Delete()
$wsOverview.Range("A4").Hyperlinks.Delete()

$wsOverview.Rows.Item(3).Delete()

$wsOverview.Range("A2").Value2 = $newMdName

$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), ($mdBase + $newMdName), [Type]::Missing, [Type]::Missing, $newMdName) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $cfgUrl, [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null

# ---------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Hyperlinks.Delete()
$wsZh.Range("C2").Hyperlinks.Delete()
$wsZh.Range("A3").Hyperlinks.Delete()
$wsZh.Range("A4").Hyperlinks.Delete()

$wsZh.Rows.Item(3).Delete()

$wsZh.Range("A2").Value2 = $newMdName
$wsZh.Range("C2").Value2 = $newZhXlfName
$wsZh.Range("D2").Value2 = $newZhDatetime

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), ($mdBase + $newMdName), [Type]::Missing, [Type]::Missing, $newMdName) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), ($zhXlfBase + $newZhXlfName), [Type]::Missing, [Type]::Missing, $newZhXlfName) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $cfgUrl, [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null

# ---------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Hyperlinks.Delete()
$wsDe.Range("C2").Hyperlinks.Delete()
$wsDe.Range("A3").Hyperlinks.Delete()
$wsDe.Range("A4").Hyperlinks.Delete()

$wsDe.Rows.Item(3).Delete()

$wsDe.Range("A2").Value2 = $newMdName
$wsDe.Range("C2").Value2 = $newDeXlfName
$wsDe.Range("D2").Value2 = $newDeDatetime

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), ($mdBase + $newMdName), [Type]::Missing, [Type]::Missing, $newMdName) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), ($deXlfBase + $newDeXlfName), [Type]::Missing, [Type]::Missing, $newDeXlfName) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $cfgUrl, [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null
